$d = $word.ActiveDocument

# 1. Remove the stray "_GoBack" bookmark that wrapped the very start of
#    the first paragraph. It is an invisible cursor-position marker left
#    over from a prior save and carries no visible content of its own.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
    Write-Output "Removed _GoBack bookmark"
}

# 2. Replace the merge-field placeholder "{{ cover_sheet_top_message }}"
#    (originally split across several runs, with spell/grammar proofing
#    marks in between) with the final, literal message text. Running
#    this through Find/Replace collapses the placeholder down to a
#    single run that keeps the bold/26pt formatting of the first run
#    that was matched.
$found = $d.Content.Find.Execute(
    "{{cover_sheet_top_message }}", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Congratulations! Your forms are ready.", 2)
if ($found) {
    Write-Output "Replaced cover_sheet_top_message placeholder"
}
